# "update excel - tareas"
#
# Renames the docker-related headers (M1/N1), turns the numeric sprint
# columns (S1:AB1) into named task labels, fills in some task-tracking
# numbers for a handful of students, and moves the current selection /
# view around a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parcial 2")
$ws.Activate()

# --- Header row: S1:AB1 used to be just sprint numbers 1..10, now they
#     carry the actual task names (write these first so new shared
#     strings land in the same order as the source workbook) -----------
$ws.Range("S1").Value = "kanban"
$ws.Range("T1").Value = "git banching"
$ws.Range("U1").Value = "issue"
$ws.Range("V1").Value = "front"
$ws.Range("W1").Value = "api"
$ws.Range("X1").Value = "token"
$ws.Range("Y1").Value = "front-api"
$ws.Range("Z1").Value = "redis"
$ws.Range("AA1").Value = "cache"
$ws.Range("AB1").Value = "tdd"

# --- Header row: rename the two "docker" task columns -----------------
$ws.Range("M1").Value = "docker maps"
$ws.Range("N1").Value = "docker api"

# --- Row 4 (Daniel Saul Chavez Garcia): fill in task scores ------------
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = 10
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 5
$ws.Range("X4").Value = 5
$ws.Range("Y4").Value = 5
$ws.Range("AA4").Value = 10

# --- Row 9 (Oscar Barbosa Aquino) ---------------------------------------
$ws.Range("L9").Value = 1

# --- Row 10 (Antonio Diaz Flores) ---------------------------------------
$ws.Range("O10").Value = 1

# --- Row 17 (Jorge Crespo Capiterucho) ----------------------------------
$ws.Range("L17").Value = 1

# --- Column widths (best effort; engine snaps to its own grid) ---------
$ws.Columns.Item(6).ColumnWidth = 9.833333333333334
$ws.Columns.Item(13).ColumnWidth = 10.666666666666666
$ws.Columns.Item(15).ColumnWidth = 16
$ws.Columns.Item(20).ColumnWidth = 10.166666666666666
$ws.Columns.Item(21).ColumnWidth = 5.833333333333333
$ws.Columns.Item(25).ColumnWidth = 11

# --- View / selection ----------------------------------------------------
$ws.Range("N9").Select()
